$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.40717887878418
$ws.Range("B1").Value = 2.548656225204468
$ws.Range("C1").Value = 2.144224882125854
$ws.Range("D1").Value = 2.229082345962524
$ws.Range("E1").Value = 2.563613176345825
